# This edit adds a new weekly price-report row for "Feria Lagunitas de
# Puerto Montt - Cilantro". The new record is inserted at row 242 (pushing
# every subsequent row down by one, so the former row 348 becomes row 349
# and the sheet's used range grows from A1:R348 to A1:R349).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 242; everything below shifts down.
$ws.Rows(242).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(242, 1).Value = 4
$ws.Cells.Item(242, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(242, 3).Value = "Los Lagos"
$ws.Cells.Item(242, 4).Value = 44845
$ws.Cells.Item(242, 5).Value = 10
$ws.Cells.Item(242, 6).Value = 100112040
$ws.Cells.Item(242, 7).Value = "Cilantro"
$ws.Cells.Item(242, 8).Value = "Sin especificar"
$ws.Cells.Item(242, 9).Value = "Primera"
$ws.Cells.Item(242, 10).Value = 180
$ws.Cells.Item(242, 11).Value = 11000
$ws.Cells.Item(242, 12).Value = 11000
$ws.Cells.Item(242, 13).Value = 11000
$ws.Cells.Item(242, 14).Value = '$/caja 36 atados'
$ws.Cells.Item(242, 15).Value = "Región Metropolitana"
$ws.Cells.Item(242, 16).Value = 306
$ws.Cells.Item(242, 17).Value = 36
$ws.Cells.Item(242, 18).Value = "Hortaliza"
